# Apply updated TPM-derived NATMI metrics to Sheet1 (Jam2-F11r LR pairs).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 7).Value = 83.40125033333334
$ws.Cells.Item(2, 8).Value = 250.203751
$ws.Cells.Item(2, 9).Value = 0.9428346765536562
$ws.Cells.Item(2, 10).Value = 0.9428346765536562
$ws.Cells.Item(2, 13).Value = 42.31746133333333
$ws.Cells.Item(2, 14).Value = 126.952384
$ws.Cells.Item(2, 15).Value = 0.6904142182914543
$ws.Cells.Item(2, 16).Value = 0.6904142182914543
$ws.Cells.Item(2, 17).Value = 3529.329186132487
$ws.Cells.Item(2, 18).Value = 31763.96267519238
$ws.Cells.Item(2, 19).Value = 0.6509464661908687
$ws.Cells.Item(2, 20).Value = 0.6509464661908687

# Row 3
$ws.Cells.Item(3, 7).Value = 83.40125033333334
$ws.Cells.Item(3, 8).Value = 250.203751
$ws.Cells.Item(3, 9).Value = 0.9428346765536562
$ws.Cells.Item(3, 10).Value = 0.9428346765536562
$ws.Cells.Item(3, 13).Value = 11.08476666666667
$ws.Cells.Item(3, 14).Value = 33.2543
$ws.Cells.Item(3, 15).Value = 0.1808492350906109
$ws.Cells.Item(3, 16).Value = 0.1808492350906109
$ws.Cells.Item(3, 17).Value = 924.4833996532556
$ws.Cells.Item(3, 18).Value = 8320.3505968793
$ws.Cells.Item(3, 19).Value = 0.1705109300716323
$ws.Cells.Item(3, 20).Value = 0.1705109300716323

# Row 4
$ws.Cells.Item(4, 7).Value = 83.40125033333334
$ws.Cells.Item(4, 8).Value = 250.203751
$ws.Cells.Item(4, 9).Value = 0.9428346765536562
$ws.Cells.Item(4, 10).Value = 0.9428346765536562
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 0.140061
$ws.Cells.Item(4, 14).Value = 0.420183
$ws.Cells.Item(4, 15).Value = 0.002285111223152439
$ws.Cells.Item(4, 16).Value = 0.002285111223152439
$ws.Cells.Item(4, 17).Value = 11.681262522937
$ws.Cells.Item(4, 18).Value = 105.131362706433
$ws.Cells.Item(4, 19).Value = 0.00215448210097006
$ws.Cells.Item(4, 20).Value = 0.00215448210097006

# Row 5
$ws.Cells.Item(5, 7).Value = 83.40125033333334
$ws.Cells.Item(5, 8).Value = 250.203751
$ws.Cells.Item(5, 9).Value = 0.9428346765536562
$ws.Cells.Item(5, 10).Value = 0.9428346765536562
$ws.Cells.Item(5, 13).Value = 7.750570000000001
$ws.Cells.Item(5, 14).Value = 23.25171
$ws.Cells.Item(5, 15).Value = 0.1264514353947823
$ws.Cells.Item(5, 16).Value = 0.1264514353947823
$ws.Cells.Item(5, 17).Value = 646.4072287960234
$ws.Cells.Item(5, 18).Value = 5817.665059164211
$ws.Cells.Item(5, 19).Value = 0.1192227981901851
$ws.Cells.Item(5, 20).Value = 0.1192227981901851

# Row 6
$ws.Cells.Item(6, 9).Value = 0.03022505171551549
$ws.Cells.Item(6, 10).Value = 0.03022505171551549
$ws.Cells.Item(6, 13).Value = 42.31746133333333
$ws.Cells.Item(6, 14).Value = 126.952384
$ws.Cells.Item(6, 15).Value = 0.6904142182914543
$ws.Cells.Item(6, 16).Value = 0.6904142182914543
$ws.Cells.Item(6, 17).Value = 113.1419535414827
$ws.Cells.Item(6, 18).Value = 1018.277581873344
$ws.Cells.Item(6, 19).Value = 0.02086780545298641
$ws.Cells.Item(6, 20).Value = 0.02086780545298641

# Row 7
$ws.Cells.Item(7, 9).Value = 0.03022505171551549
$ws.Cells.Item(7, 10).Value = 0.03022505171551549
$ws.Cells.Item(7, 13).Value = 11.08476666666667
$ws.Cells.Item(7, 14).Value = 33.2543
$ws.Cells.Item(7, 15).Value = 0.1808492350906109
$ws.Cells.Item(7, 16).Value = 0.1808492350906109
$ws.Cells.Item(7, 17).Value = 29.63675314403334
$ws.Cells.Item(7, 18).Value = 266.7307782963
$ws.Cells.Item(7, 19).Value = 0.005466177483325133
$ws.Cells.Item(7, 20).Value = 0.005466177483325133

# Row 8
$ws.Cells.Item(8, 9).Value = 0.03022505171551549
$ws.Cells.Item(8, 10).Value = 0.03022505171551549
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 0.140061
$ws.Cells.Item(8, 14).Value = 0.420183
$ws.Cells.Item(8, 15).Value = 0.002285111223152439
$ws.Cells.Item(8, 16).Value = 0.002285111223152439
$ws.Cells.Item(8, 17).Value = 0.374473672467
$ws.Cells.Item(8, 18).Value = 3.370263052203
$ws.Cells.Item(8, 19).Value = 0.00006906760489548732
$ws.Cells.Item(8, 20).Value = 0.00006906760489548732

# Row 9
$ws.Cells.Item(9, 9).Value = 0.03022505171551549
$ws.Cells.Item(9, 10).Value = 0.03022505171551549
$ws.Cells.Item(9, 13).Value = 7.750570000000001
$ws.Cells.Item(9, 14).Value = 23.25171
$ws.Cells.Item(9, 15).Value = 0.1264514353947823
$ws.Cells.Item(9, 16).Value = 0.1264514353947823
$ws.Cells.Item(9, 17).Value = 20.72228822879
$ws.Cells.Item(9, 18).Value = 186.50059405911
$ws.Cells.Item(9, 19).Value = 0.00382200117430846
$ws.Cells.Item(9, 20).Value = 0.00382200117430846

# Row 10
$ws.Cells.Item(10, 5).Value = 2
$ws.Cells.Item(10, 6).Value = 0.6666666666666666
$ws.Cells.Item(10, 7).Value = 0.2062613333333333
$ws.Cells.Item(10, 8).Value = 0.618784
$ws.Cells.Item(10, 9).Value = 0.002331743669568637
$ws.Cells.Item(10, 10).Value = 0.002331743669568637
$ws.Cells.Item(10, 13).Value = 42.31746133333333
$ws.Cells.Item(10, 14).Value = 126.952384
$ws.Cells.Item(10, 15).Value = 0.6904142182914543
$ws.Cells.Item(10, 16).Value = 0.6904142182914543
$ws.Cells.Item(10, 17).Value = 8.728455997895111
$ws.Cells.Item(10, 18).Value = 78.556103981056
$ws.Cells.Item(10, 19).Value = 0.001609868982881278
$ws.Cells.Item(10, 20).Value = 0.001609868982881278

# Row 11
$ws.Cells.Item(11, 5).Value = 2
$ws.Cells.Item(11, 6).Value = 0.6666666666666666
$ws.Cells.Item(11, 7).Value = 0.2062613333333333
$ws.Cells.Item(11, 8).Value = 0.618784
$ws.Cells.Item(11, 9).Value = 0.002331743669568637
$ws.Cells.Item(11, 10).Value = 0.002331743669568637
$ws.Cells.Item(11, 13).Value = 11.08476666666667
$ws.Cells.Item(11, 14).Value = 33.2543
$ws.Cells.Item(11, 15).Value = 0.1808492350906109
$ws.Cells.Item(11, 16).Value = 0.1808492350906109
$ws.Cells.Item(11, 17).Value = 2.286358752355556
$ws.Cells.Item(11, 18).Value = 20.5772287712
$ws.Cells.Item(11, 19).Value = 0.0004216940590688622
$ws.Cells.Item(11, 20).Value = 0.0004216940590688622

# Row 12
$ws.Cells.Item(12, 5).Value = 2
$ws.Cells.Item(12, 6).Value = 0.6666666666666666
$ws.Cells.Item(12, 7).Value = 0.2062613333333333
$ws.Cells.Item(12, 8).Value = 0.618784
$ws.Cells.Item(12, 9).Value = 0.002331743669568637
$ws.Cells.Item(12, 10).Value = 0.002331743669568637
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 12).Value = 1
$ws.Cells.Item(12, 13).Value = 0.140061
$ws.Cells.Item(12, 14).Value = 0.420183
$ws.Cells.Item(12, 15).Value = 0.002285111223152439
$ws.Cells.Item(12, 16).Value = 0.002285111223152439
$ws.Cells.Item(12, 17).Value = 0.028889168608
$ws.Cells.Item(12, 18).Value = 0.260002517472
$ws.Cells.Item(12, 19).Value = 0.000005328293628845945
$ws.Cells.Item(12, 20).Value = 0.000005328293628845945

# Row 13
$ws.Cells.Item(13, 5).Value = 2
$ws.Cells.Item(13, 6).Value = 0.6666666666666666
$ws.Cells.Item(13, 7).Value = 0.2062613333333333
$ws.Cells.Item(13, 8).Value = 0.618784
$ws.Cells.Item(13, 9).Value = 0.002331743669568637
$ws.Cells.Item(13, 10).Value = 0.002331743669568637
$ws.Cells.Item(13, 13).Value = 7.750570000000001
$ws.Cells.Item(13, 14).Value = 23.25171
$ws.Cells.Item(13, 15).Value = 0.1264514353947823
$ws.Cells.Item(13, 16).Value = 0.1264514353947823
$ws.Cells.Item(13, 17).Value = 1.598642902293334
$ws.Cells.Item(13, 18).Value = 14.38778612064
$ws.Cells.Item(13, 19).Value = 0.0002948523339896511
$ws.Cells.Item(13, 20).Value = 0.0002948523339896511

# Row 14
$ws.Cells.Item(14, 7).Value = 1.881585
$ws.Cells.Item(14, 8).Value = 5.644755
$ws.Cells.Item(14, 9).Value = 0.02127094711161878
$ws.Cells.Item(14, 10).Value = 0.02127094711161878
$ws.Cells.Item(14, 13).Value = 42.31746133333333
$ws.Cells.Item(14, 14).Value = 126.952384
$ws.Cells.Item(14, 15).Value = 0.6904142182914543
$ws.Cells.Item(14, 16).Value = 0.6904142182914543
$ws.Cells.Item(14, 17).Value = 79.62390048288
$ws.Cells.Item(14, 18).Value = 716.6151043459199
$ws.Cells.Item(14, 19).Value = 0.01468576432238715
$ws.Cells.Item(14, 20).Value = 0.01468576432238715

# Row 15
$ws.Cells.Item(15, 7).Value = 1.881585
$ws.Cells.Item(15, 8).Value = 5.644755
$ws.Cells.Item(15, 9).Value = 0.02127094711161878
$ws.Cells.Item(15, 10).Value = 0.02127094711161878
$ws.Cells.Item(15, 13).Value = 11.08476666666667
$ws.Cells.Item(15, 14).Value = 33.2543
$ws.Cells.Item(15, 15).Value = 0.1808492350906109
$ws.Cells.Item(15, 16).Value = 0.1808492350906109
$ws.Cells.Item(15, 17).Value = 20.8569306885
$ws.Cells.Item(15, 18).Value = 187.7123761965
$ws.Cells.Item(15, 19).Value = 0.003846834514789095
$ws.Cells.Item(15, 20).Value = 0.003846834514789095

# Row 16
$ws.Cells.Item(16, 7).Value = 1.881585
$ws.Cells.Item(16, 8).Value = 5.644755
$ws.Cells.Item(16, 9).Value = 0.02127094711161878
$ws.Cells.Item(16, 10).Value = 0.02127094711161878
$ws.Cells.Item(16, 11).Value = 3
$ws.Cells.Item(16, 12).Value = 1
$ws.Cells.Item(16, 13).Value = 0.140061
$ws.Cells.Item(16, 14).Value = 0.420183
$ws.Cells.Item(16, 15).Value = 0.002285111223152439
$ws.Cells.Item(16, 16).Value = 0.002285111223152439
$ws.Cells.Item(16, 17).Value = 0.263536676685
$ws.Cells.Item(16, 18).Value = 2.371830090165
$ws.Cells.Item(16, 19).Value = 0.00004860647997184202
$ws.Cells.Item(16, 20).Value = 0.00004860647997184202

# Row 17
$ws.Cells.Item(17, 7).Value = 1.881585
$ws.Cells.Item(17, 8).Value = 5.644755
$ws.Cells.Item(17, 9).Value = 0.02127094711161878
$ws.Cells.Item(17, 10).Value = 0.02127094711161878
$ws.Cells.Item(17, 13).Value = 7.750570000000001
$ws.Cells.Item(17, 14).Value = 23.25171
$ws.Cells.Item(17, 15).Value = 0.1264514353947823
$ws.Cells.Item(17, 16).Value = 0.1264514353947823
$ws.Cells.Item(17, 17).Value = 14.58335625345
$ws.Cells.Item(17, 18).Value = 131.25020628105
$ws.Cells.Item(17, 19).Value = 0.002689741794470693
$ws.Cells.Item(17, 20).Value = 0.002689741794470693

# Row 18
$ws.Cells.Item(18, 5).Value = 3
$ws.Cells.Item(18, 6).Value = 1
$ws.Cells.Item(18, 7).Value = 0.2952356666666667
$ws.Cells.Item(18, 8).Value = 0.885707
$ws.Cells.Item(18, 9).Value = 0.003337580949640955
$ws.Cells.Item(18, 10).Value = 0.003337580949640955
$ws.Cells.Item(18, 13).Value = 42.31746133333333
$ws.Cells.Item(18, 14).Value = 126.952384
$ws.Cells.Item(18, 15).Value = 0.6904142182914543
$ws.Cells.Item(18, 16).Value = 0.6904142182914543
$ws.Cells.Item(18, 17).Value = 12.49362390838756
$ws.Cells.Item(18, 18).Value = 112.442615175488
$ws.Cells.Item(18, 19).Value = 0.00230431334233081
$ws.Cells.Item(18, 20).Value = 0.00230431334233081

# Row 19
$ws.Cells.Item(19, 5).Value = 3
$ws.Cells.Item(19, 6).Value = 1
$ws.Cells.Item(19, 7).Value = 0.2952356666666667
$ws.Cells.Item(19, 8).Value = 0.885707
$ws.Cells.Item(19, 9).Value = 0.003337580949640955
$ws.Cells.Item(19, 10).Value = 0.003337580949640955
$ws.Cells.Item(19, 13).Value = 11.08476666666667
$ws.Cells.Item(19, 14).Value = 33.2543
$ws.Cells.Item(19, 15).Value = 0.1808492350906109
$ws.Cells.Item(19, 16).Value = 0.1808492350906109
$ws.Cells.Item(19, 17).Value = 3.272618476677778
$ws.Cells.Item(19, 18).Value = 29.4535662901
$ws.Cells.Item(19, 19).Value = 0.0006035989617955615
$ws.Cells.Item(19, 20).Value = 0.0006035989617955615

# Row 20
$ws.Cells.Item(20, 5).Value = 3
$ws.Cells.Item(20, 6).Value = 1
$ws.Cells.Item(20, 7).Value = 0.2952356666666667
$ws.Cells.Item(20, 8).Value = 0.885707
$ws.Cells.Item(20, 9).Value = 0.003337580949640955
$ws.Cells.Item(20, 10).Value = 0.003337580949640955
$ws.Cells.Item(20, 11).Value = 3
$ws.Cells.Item(20, 12).Value = 1
$ws.Cells.Item(20, 13).Value = 0.140061
$ws.Cells.Item(20, 14).Value = 0.420183
$ws.Cells.Item(20, 15).Value = 0.002285111223152439
$ws.Cells.Item(20, 16).Value = 0.002285111223152439
$ws.Cells.Item(20, 17).Value = 0.041351002709
$ws.Cells.Item(20, 18).Value = 0.372159024381
$ws.Cells.Item(20, 19).Value = 0.000007626743686204322
$ws.Cells.Item(20, 20).Value = 0.000007626743686204323

# Row 21
$ws.Cells.Item(21, 5).Value = 3
$ws.Cells.Item(21, 6).Value = 1
$ws.Cells.Item(21, 7).Value = 0.2952356666666667
$ws.Cells.Item(21, 8).Value = 0.885707
$ws.Cells.Item(21, 9).Value = 0.003337580949640955
$ws.Cells.Item(21, 10).Value = 0.003337580949640955
$ws.Cells.Item(21, 13).Value = 7.750570000000001
$ws.Cells.Item(21, 14).Value = 23.25171
$ws.Cells.Item(21, 15).Value = 0.1264514353947823
$ws.Cells.Item(21, 16).Value = 0.1264514353947823
$ws.Cells.Item(21, 17).Value = 2.288244700996667
$ws.Cells.Item(21, 18).Value = 20.59420230897
$ws.Cells.Item(21, 19).Value = 0.0004220419018283793
$ws.Cells.Item(21, 20).Value = 0.0004220419018283794

Write-Output "Applied 258 cell updates"
